# Update the loading_percent values on Sheet1 for rows 2-25.
# Columns touched: B, D, E, F, G, H, I, O (in that order for each row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("B", "D", "E", "F", "G", "H", "I", "O")

$newValues = @{
    2 = @(15.46388812271717,3.332866452911979,16.7405582290129,18.35869918359322,21.97987378217738,10.62994996385521,42.00370031271751,15.51357751095988)
    3 = @(14.62529716494162,3.331906116542084,15.77358072534926,17.85750361967092,20.88400046637141,10.56398091616973,39.58439637612544,15.21194444154728)
    4 = @(14.08227189835594,3.331667210034011,15.15375536400619,17.54855801098144,20.18742174270481,10.52596340044021,38.02187358361297,15.02856969616999)
    5 = @(13.85402509196118,3.33165832644058,14.89487218464948,17.4225731203034,19.89801699846327,10.51110920613297,37.36626242888532,14.9544145422213)
    6 = @(13.81570835831321,3.331662202569909,14.85151310981644,17.40165469098864,19.84964027438417,10.50868157673871,37.25627603765865,14.94213892303543)
    7 = @(14.07922169308636,3.331666731627296,15.1502890936753,17.54685899359088,20.18354057289569,10.52576047143519,38.01310743513349,15.02756715204569)
    8 = @(15.18063999774328,3.332462701585339,16.41270293448153,18.18626524397519,21.60716863689971,10.60669544107147,41.18578900812587,15.4092615210346)
    9 = @(17.11354650569201,3.336794265068324,18.83526645892429,19.42152583551522,24.19549305282739,10.78454097625027,46.77879708202151,16.16742811790309)
    10 = @(18.39153772040988,3.341648870266802,20.51580639411464,20.30649208104569,25.95640391697248,10.92598972804872,50.48748964831373,16.72400890826519)
    11 = @(18.94157951104755,3.344216520351163,21.23749084471423,20.70207008789831,26.72430725654529,10.99247610726923,52.08534703260537,16.97584131324049)
    12 = @(19.14533481163442,3.345240104828057,21.50465745367265,20.85070168403498,27.01015662864143,11.01794236505454,52.67744494454099,17.07090818264913)
    13 = @(19.10165446781605,3.345017383736601,21.44738995832763,20.81874541673822,26.94881594057337,11.01244518410014,52.55050504930303,17.05044857741606)
    14 = @(18.95843357864365,3.344299706018385,21.25959324810564,20.71432231239457,26.74792412138484,10.99456554138558,52.13432002565632,16.98366902686251)
    15 = @(18.87011552320599,3.343866772780879,21.14376642488964,20.65020392082526,26.62422428068378,10.98365086186756,51.87770123785431,16.94272308504471)
    16 = @(18.3549581758265,3.341488264673342,20.46778435323411,20.28048541154647,25.9055371870232,10.921686192136,50.38125698588405,16.70751495472316)
    17 = @(18.03087938850338,3.340120871658005,20.04215383634661,20.05176529754235,25.45602378938278,10.88420815208812,49.44026495104949,16.56280069720605)
    18 = @(17.84153069387959,3.339368223788723,19.79331164051828,19.91956150706805,25.19436406922753,10.86285440835427,48.89063678530206,16.47944270756341)
    19 = @(17.77691535278072,3.339119213877821,19.70836458320055,19.87469300417116,25.10524197005105,10.85565975489246,48.7031053676155,16.45120130807218)
    20 = @(18.06568338928438,3.340262933356617,20.08787984548762,20.07618141573258,25.50419863862297,10.88817691407267,49.54130510666898,16.57821917134406)
    21 = @(19.00062425561376,3.344509117089792,21.31491949976721,20.74502673075174,26.80706618550222,10.99980952428328,52.25691677039583,17.00329262939363)
    22 = @(19.58522931221873,3.347582896093313,22.08120714169245,21.17530009767951,27.62971804028947,11.07444526050048,53.9560509119212,17.27933077334015)
    23 = @(19.27564066207565,3.345915169760891,21.67547575067288,20.94633133710848,27.19334216302269,11.03446353262154,53.05615300520156,17.13219820137929)
    24 = @(18.04995794057753,3.340198602956222,20.06722001980252,20.06514508651011,25.48242883184155,10.88638203454029,49.49565175058866,16.5712489692286)
    25 = @(16.61539234780995,3.335327506432339,18.17820302606687,19.09055749939812,23.51897299350288,10.73446794358163,45.33509062920486,15.96195547730644)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
